# Update the per-topic slides (Slide 2 .. Slide 11 in the deck) so that
# each slide's Title, bolded body text, and the small "source" textbox at
# the bottom all get their new copy, per the commit's restructuring of the
# Quantum_presentation.pptx deck into a numbered "Slide N: ..." outline with
# citation-style sources instead of raw URLs.
#
# NOTE: this runtime's PowerShell-style engine only reliably binds
# *positional* function arguments, so avoid `param()` / `-Name value`
# syntax here.

$p = $ppt.ActivePresentation

function Set-SlideTexts($Index, $Title, $Content, $Source) {
    $slide = $p.Slides.Item($Index)

    # Shape 1: Title placeholder
    $slide.Shapes.Item(1).TextFrame.TextRange.Paragraphs(1, 1).Runs(1).Text = $Title

    # Shape 2: bold body / content placeholder
    $slide.Shapes.Item(2).TextFrame.TextRange.Paragraphs(1, 1).Runs(1).Text = $Content

    # Shape 3: small textbox with the "- source" line (paragraph 1 is a
    # blank spacer line above it, so only touch paragraph 2's run).
    $slide.Shapes.Item(3).TextFrame.TextRange.Paragraphs(2, 1).Runs(1).Text = $Source
}

Set-SlideTexts 2 "Slide 1: Introduction to Quantum" "Description of what quantum mechanics is and its importance in modern physics." "- Bohr, N. (1913). On the Constitution of Atoms and Molecules."

Set-SlideTexts 3 "Slide 2: Quantum Superposition" "Explanation of the concept of superposition and its role in quantum theory." "- Schrodinger, E. (1935). Discussion of Probability Relations between Separated Systems."

Set-SlideTexts 4 "Slide 3: Quantum Entanglement" "Explanation of entanglement and its implications for quantum computing and communication." "- Einstein, A., Podolsky, B., & Rosen, N. (1935). Can Quantum-Mechanical Description of Physical Reality Be Considered Complete?"

Set-SlideTexts 5 "Slide 4: Quantum Tunneling" "Description of tunneling phenomenon and its applications in various fields." "- Bardeen, J. (1957). Theory of Superconductivity."

Set-SlideTexts 6 "Slide 5: Quantum Teleportation" "Explanation of teleportation using quantum principles and its potential for secure communication." "- Bennett, C. H., & Brassard, G. (1993). Quantum cryptography: Public key distribution and coin tossing."

Set-SlideTexts 7 "Slide 6: Quantum Computing" "Overview of quantum computing and its advantages over classical computing." "- Feynman, R. (1982). Simulating Physics with Computers."

Set-SlideTexts 8 "Slide 7: Quantum Cryptography" "Introduction to quantum cryptography and its role in secure communication." "- Gisin, N., Ribordy, G., Tittel, W., & Zbinden, H. (2002). Quantum cryptography."

Set-SlideTexts 9 "Slide 8: Quantum Algorithms" "Overview of key algorithms used in quantum computing such as Grover's and Shor's algorithms." "- Grover, L. K. (1996). A fast quantum mechanical algorithm for database search."

Set-SlideTexts 10 "Slide 9: Applications of Quantum Mechanics" "Explanation of real-world applications of quantum mechanics in areas like medicine, materials science, and more." "- Haroche, S., & Raimond, J. M. (2006). Quantum information processing."

Set-SlideTexts 11 "Slide 10: Future of Quantum Technology" "Discussion on the potential advancements and future developments in quantum technology." "- Monroe, C., Meekhof, D. M., King, B. E., & Jefferts, S. R. (1996). Demonstration of a fundamental quantum logic gate."
